$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2375.4285
$ws.Range("J48").Value = 1166
$ws.Range("L48").Value = 3498
$ws.Range("N48").Value = -4082

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 2375.4285
$ws.Range("J56").Value = 1166
$ws.Range("L56").Value = 3498
$ws.Range("N56").Value = -4566

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5093.625
$ws.Range("J64").Value = 7067.7144
$ws.Range("L64").Value = 7067.7144
$ws.Range("N64").Value = -7563.7144

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5093.625
$ws.Range("J67").Value = 7067.7144
$ws.Range("L67").Value = 7067.7144
$ws.Range("N67").Value = -8783.714400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2975.182
$ws.Range("I86").Value = 2564.8572
$ws.Range("K86").Value = 2564.8572
$ws.Range("M86").Value = -1441.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2975.182
$ws.Range("I89").Value = 2564.8572
$ws.Range("K89").Value = 12824.286
$ws.Range("M89").Value = -7208.286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 33500.117
$ws.Range("J112").Value = 35416
$ws.Range("L112").Value = 106248
$ws.Range("N112").Value = -108464

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6413104.5
$ws.Range("I138").Value = 1534
$ws.Range("J138").Value = 7815635.5
$ws.Range("K138").Value = 4602
$ws.Range("L138").Value = 23446906.5
$ws.Range("M138").Value = 538
$ws.Range("N138").Value = -23457186.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10080.877
$ws.Range("I32").Value = 4703.9243
$ws.Range("J32").Value = 24329.8
$ws.Range("K32").Value = 4703.9243
$ws.Range("L32").Value = 24329.8
$ws.Range("M32").Value = -4416.9243
$ws.Range("N32").Value = -24903.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6101.857
$ws.Range("I45").Value = 6816.222
$ws.Range("K45").Value = 6816.222
$ws.Range("M45").Value = -6439.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4777.0264
$ws.Range("I110").Value = 4875.067
$ws.Range("J110").Value = 4409.375
$ws.Range("K110").Value = 4875.067
$ws.Range("L110").Value = 4409.375
$ws.Range("M110").Value = -2830.067
$ws.Range("N110").Value = -8499.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H48").Value = 389684
$ws.Range("J48").Value = 389684
$ws.Range("L48").Value = 389684
$ws.Range("N48").Value = -390514

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1857.4
$ws.Range("I105").Value = 1797.1818
$ws.Range("J105").Value = 2141.2856
$ws.Range("K105").Value = 1797.1818
$ws.Range("L105").Value = 2141.2856
$ws.Range("M105").Value = -50.18180000000007
$ws.Range("N105").Value = -5635.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2669.3157
$ws.Range("I58").Value = 1559.2858
$ws.Range("K58").Value = 1559.2858
$ws.Range("M58").Value = -1356.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2548.8235
$ws.Range("I99").Value = 2191.3
$ws.Range("J99").Value = 3059.5715
$ws.Range("K99").Value = 2191.3
$ws.Range("L99").Value = 3059.5715
$ws.Range("M99").Value = -693.3000000000002
$ws.Range("N99").Value = -6055.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2548.8235
$ws.Range("I126").Value = 2191.3
$ws.Range("J126").Value = 3059.5715
$ws.Range("K126").Value = 6573.900000000001
$ws.Range("L126").Value = 9178.7145
$ws.Range("M126").Value = -4103.900000000001
$ws.Range("N126").Value = -14118.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4095.04
$ws.Range("I132").Value = 4075.9048
$ws.Range("K132").Value = 12227.7144
$ws.Range("M132").Value = -9697.714399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2669.3157
$ws.Range("I136").Value = 1559.2858
$ws.Range("K136").Value = 4677.857400000001
$ws.Range("M136").Value = -2127.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.814816
$ws.Range("I2").Value = 46.75
$ws.Range("J2").Value = 8.428572000000001
$ws.Range("K2").Value = 280.5
$ws.Range("L2").Value = 50.571432
$ws.Range("M2").Value = -167.5
$ws.Range("N2").Value = -276.571432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 65.166664
$ws.Range("I17").Value = 40.875
$ws.Range("J17").Value = 113.75
$ws.Range("K17").Value = 122.625
$ws.Range("L17").Value = 341.25
$ws.Range("M17").Value = 46.375
$ws.Range("N17").Value = -679.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 631.7273
$ws.Range("I34").Value = 199.66667
$ws.Range("J34").Value = 793.75
$ws.Range("K34").Value = 599.00001
$ws.Range("L34").Value = 2381.25
$ws.Range("M34").Value = -515.00001
$ws.Range("N34").Value = -2549.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 9884.538
$ws.Range("J39").Value = 9884.538
$ws.Range("L39").Value = 29653.614
$ws.Range("N39").Value = -30241.614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1774.1177
$ws.Range("J55").Value = 2206.6667
$ws.Range("L55").Value = 6620.000100000001
$ws.Range("N55").Value = -6974.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1333
$ws.Range("J92").Value = 1458
$ws.Range("L92").Value = 4374
$ws.Range("N92").Value = -6870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3846948.5
$ws.Range("I113").Value = 7143511
$ws.Range("J113").Value = 958.8333
$ws.Range("K113").Value = 21430533
$ws.Range("L113").Value = 2876.4999
$ws.Range("M113").Value = -21428363
$ws.Range("N113").Value = -7216.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 31499.766
$ws.Range("I131").Value = 167946.17
$ws.Range("J131").Value = 2261.25
$ws.Range("K131").Value = 503838.51
$ws.Range("L131").Value = 6783.75
$ws.Range("M131").Value = -498798.51
$ws.Range("N131").Value = -16863.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2803.0908
$ws.Range("I122").Value = 2373
$ws.Range("J122").Value = 3555.75
$ws.Range("K122").Value = 7119
$ws.Range("L122").Value = 10667.25
$ws.Range("M122").Value = -4669
$ws.Range("N122").Value = -15567.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 23019
$ws.Range("I126").Value = 27167.777
$ws.Range("J126").Value = 4349.5
$ws.Range("K126").Value = 81503.33099999999
$ws.Range("L126").Value = 13048.5
$ws.Range("M126").Value = -79033.33099999999
$ws.Range("N126").Value = -17988.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2662.4167
$ws.Range("I22").Value = 1908.6666
$ws.Range("J22").Value = 2913.6667
$ws.Range("K22").Value = 1908.6666
$ws.Range("L22").Value = 2913.6667
$ws.Range("M22").Value = -1613.6666
$ws.Range("N22").Value = -3503.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2662.4167
$ws.Range("I27").Value = 1908.6666
$ws.Range("J27").Value = 2913.6667
$ws.Range("K27").Value = 1908.6666
$ws.Range("L27").Value = 2913.6667
$ws.Range("M27").Value = -1801.6666
$ws.Range("N27").Value = -3127.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2741.3333
$ws.Range("I68").Value = 2689.9
$ws.Range("K68").Value = 2689.9
$ws.Range("M68").Value = -1940.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2741.3333
$ws.Range("I71").Value = 2689.9
$ws.Range("K71").Value = 13449.5
$ws.Range("M71").Value = -9705.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 23298.6
$ws.Range("J28").Value = 23298.6
$ws.Range("L28").Value = 23298.6
$ws.Range("N28").Value = -23994.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 14017.583
$ws.Range("I94").Value = 16999.666
$ws.Range("J94").Value = 13023.556
$ws.Range("K94").Value = 16999.666
$ws.Range("L94").Value = 13023.556
$ws.Range("M94").Value = -16098.666
$ws.Range("N94").Value = -14825.556

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1660.8438
$ws.Range("J122").Value = 1964.8334
$ws.Range("L122").Value = 5894.5002
$ws.Range("N122").Value = -10794.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3899.8
$ws.Range("I126").Value = 3874.75
$ws.Range("K126").Value = 11624.25
$ws.Range("M126").Value = -9154.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 104900
$ws.Range("J135").Value = 104900
$ws.Range("L135").Value = 104900
$ws.Range("N135").Value = -115040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 45523.22
$ws.Range("I136").Value = 72896.14
$ws.Range("K136").Value = 218688.42
$ws.Range("M136").Value = -216138.42
